$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new daily status block for 26/10/2021 (rows 51-55)
$ws.Range("A51").Value = "26/10/2021"
$ws.Range("B51").Value = "Syncup on git link issue"
$ws.Range("C51").Value = "editing the command on the yavta codes"

$ws.Range("B52").Value = "clarified the doubts on the commads"
$ws.Range("C52").Value = "Bugs in adb log files, trying to fix them"

$ws.Range("B53").Value = "Internal discussion with the teammate "

$ws.Range("B54").Value = "Explored more on ffmpeg"

$ws.Range("B55").Value = "explored more on v4l2"

# Update the view to reflect scrolling/selection to the newly added rows
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B55").Select()
